# Fix typo in the "Sucursal" column values: "caracol" -> "caraco"
# (the three data rows all shared the same string, so this single
# shared-string edit updates every occurrence).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "Sucursal galería caraco;kori"
$ws.Range("G3").Value = "Sucursal galería caraco;kori"
$ws.Range("G4").Value = "Sucursal galería caraco;kori"

# Move the active selection to H3 (matches the saved workbook's cursor
# position after the edit).
[void]$ws.Range("H3").Select()
